$d = $word.ActiveDocument

# 1. Merge the split date "0" + "9" + ".05.2024" runs back into a single
#    run reading "09.05.2024" (same text, same formatting - Find/Replace
#    across the existing runs causes Word to consolidate them into one
#    run that carries the original run formatting).
$d.Content.Find.Execute("09.05.2024", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "09.05.2024", 2) | Out-Null

# 2. Add three additional empty paragraphs right before the last
#    (already empty) paragraph of the document.
$lastRange = $d.Paragraphs.Last.Range
$lastRange.InsertParagraphBefore()
$lastRange.InsertParagraphBefore()
$lastRange.InsertParagraphBefore()
